$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The player table rows were reshuffled:
#  - Row 3 (Austin Reaves)  <-> Row 14 (Tyrese Maxey)
#  - Rows 12-13 (Trey Murphy III, Jalen Johnson) <-> Rows 15-16 (Kevin Durant, Jalen Duren)
# Capture the original row values before overwriting anything.
$row3  = $ws.Range("A3:C3").Value2
$row12 = $ws.Range("A12:C12").Value2
$row13 = $ws.Range("A13:C13").Value2
$row14 = $ws.Range("A14:C14").Value2
$row15 = $ws.Range("A15:C15").Value2
$row16 = $ws.Range("A16:C16").Value2

$ws.Range("A3:C3").Value2   = $row14
$ws.Range("A12:C12").Value2 = $row15
$ws.Range("A13:C13").Value2 = $row16
$ws.Range("A14:C14").Value2 = $row3
$ws.Range("A15:C15").Value2 = $row12
$ws.Range("A16:C16").Value2 = $row13
